$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Recalculated distance table (question 1 & 2 solve) ---
# Row 3 (angle 0)
$ws.Range("C3").Value = 415.692497577683
$ws.Range("D3").Value = 466.091395548792
$ws.Range("E3").Value = 516.490293519902
$ws.Range("F3").Value = 566.889191491012
$ws.Range("G3").Value = 617.288089462121
$ws.Range("H3").Value = 667.686987433231
$ws.Range("I3").Value = 718.085885404341
$ws.Range("J3").Value = 768.484783375451

# Row 4 (angle 45)
$ws.Range("C4").Value = 415.692345731752
$ws.Range("D4").Value = 451.329735233671
$ws.Range("E4").Value = 486.967124735589
$ws.Range("F4").Value = 522.604514237507
$ws.Range("G4").Value = 558.241903739425
$ws.Range("H4").Value = 593.879293241343
$ws.Range("I4").Value = 629.516682743261
$ws.Range("J4").Value = 665.154072245179

# Row 6 (angle 135)
$ws.Range("C6").Value = 415.692345731752
$ws.Range("D6").Value = 380.054956229834
$ws.Range("E6").Value = 344.417566727916
$ws.Range("F6").Value = 308.780177225998
$ws.Range("G6").Value = 273.14278772408
$ws.Range("H6").Value = 237.505398222162
$ws.Range("I6").Value = 201.868008720244
$ws.Range("J6").Value = 166.230619218326

# Row 7 (angle 180)
$ws.Range("C7").Value = 415.692497577683
$ws.Range("D7").Value = 365.293599606573
$ws.Range("E7").Value = 314.894701635463
$ws.Range("F7").Value = 264.495803664353
$ws.Range("G7").Value = 214.096905693244
$ws.Range("H7").Value = 163.698007722134
$ws.Range("I7").Value = 113.299109751024
$ws.Range("J7").Value = 62.9002117799145

# Row 8 (angle 225)
$ws.Range("C8").Value = 415.692345731752
$ws.Range("D8").Value = 380.054956229834
$ws.Range("E8").Value = 344.417566727916
$ws.Range("F8").Value = 308.780177225998
$ws.Range("G8").Value = 273.14278772408
$ws.Range("H8").Value = 237.505398222162
$ws.Range("I8").Value = 201.868008720244
$ws.Range("J8").Value = 166.230619218326

# Row 10 (angle 315)
$ws.Range("C10").Value = 415.692345731752
$ws.Range("D10").Value = 451.329735233671
$ws.Range("E10").Value = 486.967124735589
$ws.Range("F10").Value = 522.604514237507
$ws.Range("G10").Value = 558.241903739425
$ws.Range("H10").Value = 593.879293241343
$ws.Range("I10").Value = 629.516682743261
$ws.Range("J10").Value = 665.154072245179

# --- New trailing helper rows (15, 16), matching the existing table style ---
$ws.Range("C3").Copy() | Out-Null
$ws.Range("I15:J15").PasteSpecial(-4122) | Out-Null
$ws.Range("G16:H16").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Rows.Item(15).RowHeight = 15.6
$ws.Rows.Item(16).RowHeight = 15.6

# --- View / selection state ---
$ws.Range("M8").Select() | Out-Null
